$d = $word.ActiveDocument

# 1. Update experience years in the PROFESSIONAL SUMMARY paragraph
$d.Content.Find.Execute(
    "Results-driven Marketing & Data Analytics Professional with 21 years of experience",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Results-driven Marketing & Data Analytics Professional with 15+ years of experience",
    2)

# 2. Enhance the FLEEM bullet under Progressive Change Campaign Committee / RESEARCH DIRECTOR
$d.Content.Find.Execute(
    "Engineered FLEEM web application using Twilio's API to make thousands of simultaneous phone calls for IVR polls",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Conceived, architected, and engineered FLEEM web application using Twilio API handling tens of thousands of calls using emulated predictive dialer for regulated political surveys",
    2)

# 3. Add a new bullet point after the "Developed innovative approaches..." paragraph
#    under the Lake Research Partners / PROGRAMMER role.
$targetIdx = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -like "*Developed innovative approaches to visualizing demographic and market data, enhancing clients' understanding of research findings*") {
        $targetIdx = $i
        break
    }
}

if ($targetIdx -gt 0) {
    $p = $d.Paragraphs.Item($targetIdx)
    $p.Range.InsertParagraphAfter()
    $newPara = $d.Paragraphs.Item($targetIdx + 1)
    $newPara.Range.InsertAfter("• Trained staff on building Python tooling for report generation and analysis")
}
